# Refresh cryptocurrency price and volume data (GitHub Actions symbol list update)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'277.22"
$ws.Range("D2").Style = "Normal"
$ws.Range("D3").Value = "'27.33"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'0.44%"
$ws.Range("E3").Style = "Normal"
$ws.Range("D4").Value = "'4.840"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'1.80%"
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'0.06323"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'0.34%"
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'7.023"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'1.41%"
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = "'1.316"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'-2.57%"
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = "'0.8875"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'1.23%"
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = "'0.1512"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'0.02%"
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'0.05396"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'7.33%"
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'0.07451"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'-0.88%"
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'0.02873"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'-3.31%"
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'0.08942"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'-0.65%"
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'0.001563"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'-0.45%"
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'0.0006347"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'-0.16%"
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'0.006033"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'0.12%"
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'3.472"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'0.74%"
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'3.298"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'0.03%"
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "'2.232"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'-2.29%"
$ws.Range("E19").Style = "Normal"
$ws.Range("D21").Value = "'0.1349"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'0.99%"
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'3.905"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'-0.29%"
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'0.1507"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'9.19%"
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'0.04389"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'-0.24%"
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'0.001178"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'0.43%"
$ws.Range("E25").Style = "Normal"
$ws.Range("E26").Value = "'10.73%"
$ws.Range("E26").Style = "Normal"
$ws.Range("E28").Value = "'-1.64%"
$ws.Range("E28").Style = "Normal"
$ws.Range("E29").Value = "'-14.83%"
$ws.Range("E29").Style = "Normal"
$ws.Range("D40").Value = "'0.03996"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'-2.57%"
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'0.006689"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'-2.11%"
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'0.1406"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'20.02%"
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'0.002140"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'0.49%"
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'0.01159"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'0.63%"
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'0.00005540"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'7.39%"
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = "'1.561"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'4.79%"
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = "'0.01850"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'-19.54%"
$ws.Range("E47").Style = "Normal"
